$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear B5 (age value 22 removed)
$ws.Range("B5").ClearContents()

# E19: change from numeric 1 to text "l"
$ws.Range("E19").Value = "l"

# E29: clear its text content (was "3 ", now empty) but keep the cell's
# existing quote-prefix style
$ws.Range("E29").ClearContents()

# Update the selected cell / active cell to G28
$ws.Range("G28").Select()
